# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns AD:AF are appended after the existing AC ("Unnamed: 28")
# column, with AD1:AF1 holding headers styled like the rest of row 1,
# and AD2:AF59 holding the same W/L/T values (69 / 93 / 0) for every
# player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/border/centered
# style used by the rest of row 1, then overwrite their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-59) shares the same team record.
$wins = 69
$losses = 93
$ties = 0

for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
